# Rename the header cell "Name" -> "Tag" and mark the header row (A1:A2)
# with an explicit cell style, matching the updated buffer layout used by
# the refactored config/plasmid-feature loader. Also move the active
# selection to A2 (first data row) to reflect the new header position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 held the shared string "Name" -> retitle it "Tag"
$ws.Range("A1").Value = "Tag"

# Apply the (Normal) cell style explicitly to the header + first data
# row so both carry their own formatting record instead of sharing the
# sheet's implicit default style.
$ws.Range("A1:A2").Style = "Normal"

# Move the frozen-pane selection up from A3 to A2.
[void]$ws.Range("A2").Select()
